# Update the cryptocurrency price/volume table to the latest scrape.
# Rows 2-23 keep their coin/link but refresh Price (D) and Volume(1h) (E).
# A new "LEO" row is inserted at row 24, shifting rows 24-45 down by one
# (row 46 "OKB" is unaffected), and the former last row ("Mantle") drops off
# the fixed 50-row table while rows 47-51 shift down by one more.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.872.95"
$ws.Range("E2").Value = "  +1.82%  "

$ws.Range("D3").Value = "2.811.95"
$ws.Range("E3").Value = "  +7.44%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'187.71"
$ws.Range("E5").Value = "  +1.74%  "

$ws.Range("D6").Value = "'594.05"
$ws.Range("E6").Value = "  +2.36%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.548"
$ws.Range("E8").Value = "  +2.98%  "

$ws.Range("E9").Value = "  -3.44%  "

$ws.Range("D10").Value = "2.812.28"
$ws.Range("E10").Value = "  +7.46%  "

$ws.Range("E11").Value = "  -1.26%  "

$ws.Range("D12").Value = "'0.371"
$ws.Range("E12").Value = "  +3.68%  "

$ws.Range("D13").Value = "'4.86"
$ws.Range("E13").Value = "  +2.54%  "

$ws.Range("D14").Value = "3.331.96"
$ws.Range("E14").Value = "  +7.58%  "

$ws.Range("D15").Value = "74.749.08"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").Value = "'26.81"
$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").Value = "2.817.76"
$ws.Range("E18").Value = "  +7.48%  "

$ws.Range("D19").Value = "'8.94"
$ws.Range("E19").Value = "  -1.58%  "

$ws.Range("E20").Value = "  +3.79%  "

$ws.Range("D21").Value = "'376.31"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("E22").Value = "  -2.25%  "

$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'6.20"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'70.63"
$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.960.23"
$ws.Range("E27").Value = "  +8.15%  "

$ws.Range("B28").Value = "NEARProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").Value = "'4.15"
$ws.Range("E28").Value = "  +0.65%  "

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'9.68"
$ws.Range("E29").Value = "  +3.42%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000104"
$ws.Range("E30").Value = "  +11.48%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.40"
$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'511.66"
$ws.Range("E33").Value = "  -2.43%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'7.70"
$ws.Range("E34").Value = "  +1.49%  "

$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "'1.78"
$ws.Range("E35").Value = "  +3.08%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'162.09"
$ws.Range("E37").Value = "  +0.51%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'19.89"
$ws.Range("E38").Value = "  +3.91%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.118"
$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "'19.37"
$ws.Range("E40").Value = "  +0.61%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'186.34"
$ws.Range("E41").Value = "  +15.77%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "'5.00"
$ws.Range("E43").Value = "  +2.51%  "

$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.338"
$ws.Range("E44").Value = "  +3.70%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.66"
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").Value = "'40.01"
$ws.Range("E46").Value = "  +2.77%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = "  +2.30%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.32"
$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0850"
$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'0.569"
$ws.Range("E50").Value = "  +7.99%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'3.69"
$ws.Range("E51").Value = "  +2.29%  "
